# registration works up to response validation
#
# Applies the diff:
#  - Title paragraph: mark "IDESGidp" and "IdentityProvider" as separate
#    spell-check exceptions (w:proofErr spellStart/spellEnd) with the
#    trailing space split into its own run.
#  - New paragraph after the title: "Test u2f at " + hyperlink to the
#    Yubico U2F demo.
#  - New paragraph after "Helped with ...": "Good explication of web auth
#    with x.509 here " + hyperlink to the x509-webauth wiki.
#  - The _GoBack bookmark moves out of the "Helped with" paragraph into
#    its own (new) empty paragraph, right after the new x.509 paragraph.

$d = $word.ActiveDocument

function New-HyperlinkXml([string]$rid, [string]$display) {
    $esc = $display -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    return "<w:hyperlink r:id='$rid' w:history='1'><w:r><w:t>$esc</w:t></w:r></w:hyperlink>"
}

function Set-ParagraphXml($paragraph, [string]$bodyInnerXml) {
    $xml = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>
<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>
<pkg:xmlData>
<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:r='http://schemas.openxmlformats.org/officeDocument/2006/relationships'>
<w:body>
$bodyInnerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $paragraph.Range.InsertXML($xml) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Title paragraph: split "IDESGidp " into its own spell-checked word
#    plus a standalone space run, matching the existing pattern already
#    used for "IdentityProvider".
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$titleInner = "<w:p>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>IDESGidp</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>IdentityProvider</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> commentary</w:t></w:r>" +
    "</w:p>"
Set-ParagraphXml $p1 $titleInner

Write-Host "Title now:" $d.Paragraphs.Item(1).Range.Text

# ---------------------------------------------------------------------
# 2) The paragraph right after the title was empty; turn it into the
#    "Test u2f at <link>" paragraph. Mint the hyperlink relationship via
#    Hyperlinks.Add (placement doesn't matter - we rewrite the paragraph
#    right after), then rebuild the paragraph referencing that relation-
#    ship id explicitly, and finally restore the Hyperlink character
#    style on the link run (direct-XML inserts don't carry rStyle).
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$null = $d.Hyperlinks.Add($p2.Range, "https://demo.yubico.com/u2f")
$yubicoRid = "rId8"

$p2 = $d.Paragraphs.Item(2)
$u2fInner = "<w:p>" +
    "<w:r><w:t xml:space='preserve'>Test u2f at </w:t></w:r>" +
    (New-HyperlinkXml $yubicoRid "https://demo.yubico.com/u2f") +
    "</w:p>"
Set-ParagraphXml $p2 $u2fInner

for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hl = $d.Hyperlinks.Item($i)
    if ($hl.Address -eq "https://demo.yubico.com/u2f") {
        $hl.Range.Style = "Hyperlink"
    }
}

Write-Host "New u2f paragraph:" $d.Paragraphs.Item(2).Range.Text

# ---------------------------------------------------------------------
# 3) The "Helped with ..." paragraph currently also carries the
#    _GoBack bookmark. Rebuild it without the bookmark (reusing its
#    existing hyperlink relationship id so the link keeps working),
#    then restore the Hyperlink style on its run.
# ---------------------------------------------------------------------
$helpedIndex = 4
$pHelped = $d.Paragraphs.Item($helpedIndex)
Write-Host "Paragraph $($helpedIndex) before:" $pHelped.Range.Text

$scottAddress = $null
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hl = $d.Hyperlinks.Item($i)
    if ($hl.Range.Start -ge $pHelped.Range.Start -and $hl.Range.Start -lt $pHelped.Range.End) {
        $scottAddress = $hl.Address
        $scottRid = "rId5"
    }
}

$helpedInner = "<w:p>" +
    "<w:r><w:t xml:space='preserve'>Helped with </w:t></w:r>" +
    (New-HyperlinkXml $scottRid "Scott Brady on IdSvr4") +
    "</w:p>"
Set-ParagraphXml $pHelped $helpedInner

for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hl = $d.Hyperlinks.Item($i)
    if ($hl.Address -eq $scottAddress) {
        $hl.Range.Style = "Hyperlink"
    }
}

Write-Host "Paragraph $($helpedIndex) after:" $d.Paragraphs.Item($helpedIndex).Range.Text

# ---------------------------------------------------------------------
# 4) Insert two brand-new paragraphs right after the "Helped with ..."
#    paragraph:
#      a) "Good explication of web auth with x.509 here <link>"
#      b) an otherwise-empty paragraph holding only the _GoBack bookmark
#    Build them via a trailing empty paragraph inserted with
#    InsertParagraphAfter, then flesh each one out with InsertXML.
# ---------------------------------------------------------------------
$pHelped = $d.Paragraphs.Item($helpedIndex)
$pHelped.Range.InsertParagraphAfter()
$pGood = $d.Paragraphs.Item($helpedIndex + 1)

$null = $d.Hyperlinks.Add($pGood.Range, "https://github.com/martinpaljak/x509-webauth/wiki/WebAuth")
$x509Rid = "rId9"

$pGood = $d.Paragraphs.Item($helpedIndex + 1)
$goodInner = "<w:p>" +
    "<w:r><w:t xml:space='preserve'>Good explication of web </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>auth</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> with x.509 here </w:t></w:r>" +
    (New-HyperlinkXml $x509Rid "https://github.com/martinpaljak/x509-webauth/wiki/WebAuth") +
    "</w:p>"
Set-ParagraphXml $pGood $goodInner

for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $hl = $d.Hyperlinks.Item($i)
    if ($hl.Address -eq "https://github.com/martinpaljak/x509-webauth/wiki/WebAuth") {
        $hl.Range.Style = "Hyperlink"
    }
}

Write-Host "Good explication paragraph:" $d.Paragraphs.Item($helpedIndex + 1).Range.Text

# ---------------------------------------------------------------------
# 5) New trailing paragraph that only holds the _GoBack bookmark.
# ---------------------------------------------------------------------
$pGood = $d.Paragraphs.Item($helpedIndex + 1)
$pGood.Range.InsertParagraphAfter()
$pBookmark = $d.Paragraphs.Item($helpedIndex + 2)

$bookmarkInner = "<w:p><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
Set-ParagraphXml $pBookmark $bookmarkInner

Write-Host "Bookmark paragraph text (should be empty):" $d.Paragraphs.Item($helpedIndex + 2).Range.Text
Write-Host "Total paragraphs:" $d.Paragraphs.Count
